# Apply cryptocurrency price/volume updates to match the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holding plain decimal-looking numbers must be forced to remain
# text (matching the original inline-string cell type) instead of being
# auto-converted to numeric values by Excel.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Now write the updated values.
$ws.Range("D2").Value = "40.126.59"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "2.240.20"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "293.29"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "87.52"
$ws.Range("E6").Value = "  +4.86%  "
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "31.46"
$ws.Range("E10").Value = "  +7.59%  "
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "47.00"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "2.588.24"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "14.15"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "2.237.03"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").Value = "40.089.17"
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").Value = "11.26"
$ws.Range("E21").Value = "  +8.01%  "
$ws.Range("D22").Value = "5.85"
$ws.Range("E22").Value = "  +2.21%  "
$ws.Range("D23").Value = "65.61"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").Value = "236.51"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").Value = "22.98"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("D30").Value = "9.34"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("D31").Value = "33.50"
$ws.Range("E31").Value = "  +4.39%  "
$ws.Range("D32").Value = "151.97"
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("D35").Value = "0.0723"
$ws.Range("E35").Value = "  +3.97%  "
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("D37").Value = "16.47"
$ws.Range("E37").Value = "  +7.95%  "
$ws.Range("E38").Value = "  +7.18%  "
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("D40").Value = "0.101"
$ws.Range("E40").Value = "  +4.44%  "
$ws.Range("D41").Value = "1.73"
$ws.Range("E41").Value = "  +5.94%  "
$ws.Range("E42").Value = "  +5.14%  "
$ws.Range("D43").Value = "2.077.46"
$ws.Range("E43").Value = "  +8.76%  "
$ws.Range("D44").Value = "18.64"
$ws.Range("E44").Value = "  +16.27%  "
$ws.Range("D45").Value = "0.0271"
$ws.Range("E45").Value = "  +4.49%  "
$ws.Range("E46").Value = "  +4.74%  "
$ws.Range("E47").Value = "  +11.59%  "
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "72.19"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.442.99"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").Value = "89.54"
$ws.Range("E51").Value = "  +2.62%  "
